$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1842.7142
$ws.Range("I18").Value = 1849.8334
$ws.Range("J18").Value = 1800
$ws.Range("K18").Value = 1849.8334
$ws.Range("L18").Value = 1800
$ws.Range("M18").Value = -1565.8334
$ws.Range("N18").Value = -2368
$ws.Range("H64").Value = 5997
$ws.Range("J64").Value = 5997
$ws.Range("L64").Value = 5997
$ws.Range("N64").Value = -6493
$ws.Range("H67").Value = 5997
$ws.Range("J67").Value = 5997
$ws.Range("L67").Value = 5997
$ws.Range("N67").Value = -7713
$ws.Range("H116").Value = 4741.8945
$ws.Range("I116").Value = 4689.087
$ws.Range("J116").Value = 4962.727
$ws.Range("K116").Value = 4689.087
$ws.Range("L116").Value = 4962.727
$ws.Range("M116").Value = -1247.087
$ws.Range("N116").Value = -11846.727
$ws.Range("H132").Value = 10759.647
$ws.Range("I132").Value = 8072.5625
$ws.Range("J132").Value = 53753
$ws.Range("K132").Value = 24217.6875
$ws.Range("L132").Value = 161259
$ws.Range("M132").Value = -21687.6875
$ws.Range("N132").Value = -166319
$ws.Range("H133").Value = 92249.5
$ws.Range("J133").Value = 92249.5
$ws.Range("L133").Value = 92249.5
$ws.Range("N133").Value = -102369.5
$ws.Range("H141").Value = 7763.2104
$ws.Range("J141").Value = 8996.666999999999
$ws.Range("L141").Value = 26990.001
$ws.Range("N141").Value = -37350.001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 4642.1113
$ws.Range("I36").Value = 2347.5
$ws.Range("K36").Value = 2347.5
$ws.Range("M36").Value = -2001.5
$ws.Range("H61").Value = 13495.8125
$ws.Range("I61").Value = 7226.3
$ws.Range("K61").Value = 7226.3
$ws.Range("M61").Value = -7014.3
$ws.Range("H63").Value = 3066.6667
$ws.Range("I63").Value = 2850
$ws.Range("K63").Value = 2850
$ws.Range("M63").Value = -2164
$ws.Range("H66").Value = 3066.6667
$ws.Range("I66").Value = 2850
$ws.Range("K66").Value = 14250
$ws.Range("M66").Value = -10818
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30676
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32340
$ws.Range("H97").Value = 2003.3784
$ws.Range("I97").Value = 684.4194
$ws.Range("J97").Value = 8818
$ws.Range("K97").Value = 684.4194
$ws.Range("L97").Value = 8818
$ws.Range("M97").Value = -188.4194
$ws.Range("N97").Value = -9810
$ws.Range("H102").Value = 36413.668
$ws.Range("I102").Value = 3696.6
$ws.Range("K102").Value = 3696.6
$ws.Range("M102").Value = -2074.6
$ws.Range("H132").Value = 2008008.5
$ws.Range("I132").Value = 3282.257
$ws.Range("K132").Value = 9846.771000000001
$ws.Range("M132").Value = -7316.771000000001
$ws.Range("H136").Value = 13495.8125
$ws.Range("I136").Value = 7226.3
$ws.Range("K136").Value = 21678.9
$ws.Range("M136").Value = -19128.9

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2395.48
$ws.Range("I94").Value = 1229.1177
$ws.Range("J94").Value = 4874
$ws.Range("K94").Value = 1229.1177
$ws.Range("L94").Value = 4874
$ws.Range("M94").Value = -778.1177
$ws.Range("N94").Value = -5776
$ws.Range("H134").Value = 10783.167
$ws.Range("I134").Value = 5086.76
$ws.Range("K134").Value = 15260.28
$ws.Range("M134").Value = -12725.28

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22372.2
$ws.Range("I31").Value = 12258.4
$ws.Range("K31").Value = 12258.4
$ws.Range("M31").Value = -11963.4
$ws.Range("H34").Value = 22372.2
$ws.Range("I34").Value = 12258.4
$ws.Range("K34").Value = 12258.4
$ws.Range("M34").Value = -12056.4
$ws.Range("H58").Value = 13534.444
$ws.Range("I58").Value = 8016.5835
$ws.Range("K58").Value = 8016.5835
$ws.Range("M58").Value = -7813.5835
$ws.Range("H122").Value = 8253.666999999999
$ws.Range("I122").Value = 6227.222
$ws.Range("K122").Value = 18681.666
$ws.Range("M122").Value = -16231.666
$ws.Range("H131").Value = 49999
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 49999
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 49999
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -60079
$ws.Range("H132").Value = 5977.811
$ws.Range("I132").Value = 2183.7144
$ws.Range("J132").Value = 10957.5625
$ws.Range("K132").Value = 6551.1432
$ws.Range("L132").Value = 32872.6875
$ws.Range("M132").Value = -4021.1432
$ws.Range("N132").Value = -37932.6875
$ws.Range("H136").Value = 13534.444
$ws.Range("I136").Value = 8016.5835
$ws.Range("K136").Value = 24049.7505
$ws.Range("M136").Value = -21499.7505

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1472.15
$ws.Range("J131").Value = 1498.5269
$ws.Range("L131").Value = 4495.5807
$ws.Range("N131").Value = -14575.5807
$ws.Range("H140").Value = 2857.889
$ws.Range("I140").Value = 1948.1666
$ws.Range("K140").Value = 5844.4998
$ws.Range("M140").Value = -664.4997999999996

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11490.5
$ws.Range("J80").Value = 17207.309
$ws.Range("L80").Value = 17207.309
$ws.Range("N80").Value = -19203.309
$ws.Range("H83").Value = 11490.5
$ws.Range("J83").Value = 17207.309
$ws.Range("L83").Value = 86036.54500000001
$ws.Range("N83").Value = -96020.54500000001
$ws.Range("H123").Value = 53049.6
$ws.Range("J123").Value = 53049.6
$ws.Range("L123").Value = 53049.6
$ws.Range("N123").Value = -57949.6

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6979.15
$ws.Range("I7").Value = 4954.6665
$ws.Range("J7").Value = 8635.546
$ws.Range("K7").Value = 4954.6665
$ws.Range("L7").Value = 8635.546
$ws.Range("M7").Value = -4842.6665
$ws.Range("N7").Value = -8859.546
$ws.Range("H93").Value = 12716.706
$ws.Range("I93").Value = 9500.299999999999
$ws.Range("J93").Value = 17311.572
$ws.Range("K93").Value = 9500.299999999999
$ws.Range("L93").Value = 17311.572
$ws.Range("M93").Value = -8252.299999999999
$ws.Range("N93").Value = -19807.572
$ws.Range("H126").Value = 6979.15
$ws.Range("I126").Value = 4954.6665
$ws.Range("J126").Value = 8635.546
$ws.Range("K126").Value = 14863.9995
$ws.Range("L126").Value = 25906.638
$ws.Range("M126").Value = -12393.9995
$ws.Range("N126").Value = -30846.638
$ws.Range("H136").Value = 10444.53
$ws.Range("I136").Value = 12379.667
$ws.Range("J136").Value = 8993.179
$ws.Range("K136").Value = 37139.001
$ws.Range("L136").Value = 26979.537
$ws.Range("M136").Value = -34589.001
$ws.Range("N136").Value = -32079.537

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 56665.668
$ws.Range("J98").Value = 55000
$ws.Range("L98").Value = 55000
$ws.Range("N98").Value = -60990
$ws.Range("H132").Value = 8199.973
$ws.Range("I132").Value = 2375.261
$ws.Range("K132").Value = 7125.782999999999
$ws.Range("M132").Value = -4595.782999999999
